$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update initial weight sliders (column B) and the computed optimal
# portfolio weights (columns C and D) to reflect the new slider values.

$ws.Range("B2").Value = 0.448
$ws.Range("C2").Value = 0.1295008802254504
$ws.Range("D2").Value = 0.1295008734111448

$ws.Range("C3").Value = 0.1316564287853225
$ws.Range("D3").Value = 0.1316564464375971

$ws.Range("C4").Value = 0.004124283841947814
$ws.Range("D4").Value = 0.004124274747723301

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.00001061149016206822
$ws.Range("D5").Value = 0.00001060878357728654

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0.006724954885701882
$ws.Range("D6").Value = 0.006724945978417399

$ws.Range("B7").Value = 0.346
$ws.Range("C7").Value = 0.1327159474968135
$ws.Range("D7").Value = 0.1327159446342338

$ws.Range("B8").Value = 0.204
$ws.Range("C8").Value = 0.5952668932746017
$ws.Range("D8").Value = 0.5952669060073064
